$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-14 21:48:48'
$ws.Range('I2').Value = '35.6 mm'
$ws.Range('N2').Value = '-3.3 °C 21:29 TU'
$ws.Range('E3').Value = '2026-02-14 21:48:51'
$ws.Range('L3').Value = '58.7 km/h - 227º 21:00 TU'
$ws.Range('N3').Value = '-7.6 °C 21:29 TU'
$ws.Range('O3').Value = '-5.3 °C'
$ws.Range('E4').Value = '2026-02-14 21:48:54'
$ws.Range('H4').Value = "'71%"
$ws.Range('C2').Copy()
$ws.Range('H4').PasteSpecial(-4122)
$ws.Range('J4').Value = '997.7 hPa'
$ws.Range('E5').Value = '2026-02-14 21:48:57'
$ws.Range('I5').Value = '22.3 mm'
$ws.Range('N5').Value = '-7.5 °C 21:29 TU'
$ws.Range('O5').Value = '-5.3 °C'
$ws.Range('E6').Value = '2026-02-14 21:49:00'
$ws.Range('H6').Value = "'73%"
$ws.Range('C2').Copy()
$ws.Range('H6').PasteSpecial(-4122)
$ws.Range('J6').Value = '997.8 hPa'
$ws.Range('E7').Value = '2026-02-14 21:49:02'
$ws.Range('J7').Value = '997.9 hPa'
$ws.Range('O7').Value = '13.4 °C'
$ws.Range('E8').Value = '2026-02-14 21:49:05'
$ws.Range('H8').Value = "'61%"
$ws.Range('C2').Copy()
$ws.Range('H8').PasteSpecial(-4122)
$ws.Range('J8').Value = '997.8 hPa'
$ws.Range('E9').Value = '2026-02-14 21:49:08'
$ws.Range('L9').Value = '73.4 km/h - 3º 21:21 TU'
$ws.Range('N9').Value = '9.4 °C 21:29 TU'
$ws.Range('E10').Value = '2026-02-14 21:49:10'
$ws.Range('H10').Value = "'75%"
$ws.Range('C2').Copy()
$ws.Range('H10').PasteSpecial(-4122)
$ws.Range('O10').Value = '10.0 °C'
$ws.Range('E11').Value = '2026-02-14 21:49:13'
$ws.Range('H11').Value = "'59%"
$ws.Range('C2').Copy()
$ws.Range('H11').PasteSpecial(-4122)
$ws.Range('E12').Value = '2026-02-14 21:49:16'
$ws.Range('N12').Value = '9.9 °C 21:27 TU'
$ws.Range('E13').Value = '2026-02-14 21:49:18'
$ws.Range('J13').Value = '1000.6 hPa'
$ws.Range('E14').Value = '2026-02-14 21:49:21'
$ws.Range('N14').Value = '7.8 °C 21:23 TU'
$ws.Range('O14').Value = '13.4 °C'
$ws.Range('E15').Value = '2026-02-14 21:49:23'
$ws.Range('N15').Value = '9.1 °C 21:29 TU'
$ws.Range('O15').Value = '11.2 °C'
$ws.Range('E16').Value = '2026-02-14 21:49:25'
$ws.Range('I16').Value = '8.7 mm'
$ws.Range('N16').Value = '-9.0 °C 21:09 TU'
$ws.Range('O16').Value = '-6.2 °C'
$ws.Range('E17').Value = '2026-02-14 21:49:28'
$ws.Range('L17').Value = '65.2 km/h - 53º 21:26 TU'
$ws.Range('N17').Value = '-1.2 °C 21:26 TU'
$ws.Range('E18').Value = '2026-02-14 21:49:31'
$ws.Range('H18').Value = "'74%"
$ws.Range('C2').Copy()
$ws.Range('H18').PasteSpecial(-4122)
$ws.Range('J18').Value = '998.0 hPa'
$ws.Range('O18').Value = '10.7 °C'
$ws.Range('E19').Value = '2026-02-14 21:49:34'
$ws.Range('H19').Value = "'75%"
$ws.Range('C2').Copy()
$ws.Range('H19').PasteSpecial(-4122)
$ws.Range('E20').Value = '2026-02-14 21:49:37'
$ws.Range('I20').Value = '4.8 mm'
$ws.Range('N20').Value = '-8.5 °C 21:29 TU'
$ws.Range('O20').Value = '-5.5 °C'
$ws.Range('E21').Value = '2026-02-14 21:49:40'
$ws.Range('H21').Value = "'68%"
$ws.Range('C2').Copy()
$ws.Range('H21').PasteSpecial(-4122)
$ws.Range('J21').Value = '1000.5 hPa'
$ws.Range('L21').Value = '85.3 km/h - 359º 21:08 TU'
$ws.Range('E22').Value = '2026-02-14 21:49:42'
$ws.Range('N22').Value = '-9.5 °C 21:00 TU'
$ws.Range('E23').Value = '2026-02-14 21:49:45'
$ws.Range('I23').Value = '40.3 mm'
$ws.Range('N23').Value = '-8.4 °C 21:04 TU'
$ws.Range('E24').Value = '2026-02-14 21:49:48'
$ws.Range('J24').Value = '1002.1 hPa'
$ws.Range('O24').Value = '9.4 °C'
$ws.Range('E25').Value = '2026-02-14 21:49:51'
$ws.Range('I25').Value = '17.8 mm'
$ws.Range('N25').Value = '-7.7 °C 21:24 TU'
$ws.Range('E26').Value = '2026-02-14 21:49:54'
$ws.Range('E27').Value = '2026-02-14 21:49:57'
$ws.Range('E28').Value = '2026-02-14 21:49:59'
$ws.Range('H28').Value = "'65%"
$ws.Range('C2').Copy()
$ws.Range('H28').PasteSpecial(-4122)
$ws.Range('J28').Value = '997.7 hPa'
$ws.Range('E29').Value = '2026-02-14 21:50:02'
$ws.Range('E30').Value = '2026-02-14 21:50:05'
$ws.Range('J30').Value = '997.7 hPa'
$ws.Range('O30').Value = '11.5 °C'
$ws.Range('E31').Value = '2026-02-14 21:50:08'
$ws.Range('J31').Value = '996.8 hPa'
$ws.Range('N31').Value = '7.9 °C 21:13 TU'
$ws.Range('E32').Value = '2026-02-14 21:50:11'
$ws.Range('O32').Value = '4.2 °C'
$ws.Range('E33').Value = '2026-02-14 21:50:13'
$ws.Range('H33').Value = "'62%"
$ws.Range('C2').Copy()
$ws.Range('H33').PasteSpecial(-4122)
$ws.Range('J33').Value = '1000.0 hPa'
$ws.Range('K33').Value = '2.6 MJ/m2'
$ws.Range('E34').Value = '2026-02-14 21:50:16'
$ws.Range('I34').Value = '3.9 mm'
$ws.Range('N34').Value = '-5.3 °C 21:27 TU'
$ws.Range('E35').Value = '2026-02-14 21:50:19'
$ws.Range('H35').Value = "'85%"
$ws.Range('C2').Copy()
$ws.Range('H35').PasteSpecial(-4122)
$ws.Range('J35').Value = '1004.6 hPa'
$ws.Range('N35').Value = '1.4 °C 21:15 TU'
$ws.Range('O35').Value = '2.6 °C'
$ws.Range('E36').Value = '2026-02-14 21:50:22'
$ws.Range('J36').Value = '998.4 hPa'
$ws.Range('N36').Value = '10.2 °C 21:24 TU'
$ws.Range('O36').Value = '11.9 °C'
$ws.Range('E37').Value = '2026-02-14 21:50:25'
$ws.Range('H37').Value = "'63%"
$ws.Range('C2').Copy()
$ws.Range('H37').PasteSpecial(-4122)
$ws.Range('J37').Value = '998.7 hPa'
$ws.Range('L37').Value = '68.4 km/h - 349º 21:22 TU'
$ws.Range('E38').Value = '2026-02-14 21:50:27'
$ws.Range('H38').Value = "'80%"
$ws.Range('C2').Copy()
$ws.Range('H38').PasteSpecial(-4122)
$ws.Range('E39').Value = '2026-02-14 21:50:30'
$ws.Range('I39').Value = '14.2 mm'
$ws.Range('N39').Value = '-8.6 °C 21:26 TU'
$ws.Range('O39').Value = '-5.9 °C'
$ws.Range('E40').Value = '2026-02-14 21:50:33'
$ws.Range('H40').Value = "'65%"
$ws.Range('C2').Copy()
$ws.Range('H40').PasteSpecial(-4122)
$ws.Range('J40').Value = '1001.2 hPa'
$ws.Range('E41').Value = '2026-02-14 21:50:35'
$ws.Range('J41').Value = '999.8 hPa'
$ws.Range('E42').Value = '2026-02-14 21:50:38'
$ws.Range('H42').Value = "'61%"
$ws.Range('C2').Copy()
$ws.Range('H42').PasteSpecial(-4122)
$ws.Range('E43').Value = '2026-02-14 21:50:41'
$ws.Range('H43').Value = "'64%"
$ws.Range('C2').Copy()
$ws.Range('H43').PasteSpecial(-4122)
$ws.Range('E44').Value = '2026-02-14 21:50:44'
$ws.Range('N44').Value = '-8.2 °C 21:28 TU'
$ws.Range('O44').Value = '-5.5 °C'
$ws.Range('E45').Value = '2026-02-14 21:50:46'
$ws.Range('G45').Value = '3 cm'
$ws.Range('I45').Value = '14.4 mm'
$ws.Range('J45').Value = '1007.1 hPa'
$ws.Range('N45').Value = '0.4 °C 21:14 TU'
$ws.Range('E46').Value = '2026-02-14 21:50:49'
$ws.Range('J46').Value = '1003.2 hPa'
$excel.CutCopyMode = $false
